# Rawdata/metadataframsteg.xlsx
#
# The author added a new task ("Distance sampling-analys på ripdata") to the
# progress tracker sheet. In the sheet, this is a brand-new row inserted at
# row 6 (pushing every following row down by one), with the "ej påbörjat"
# (not started) status in column B, styled the same (red font) as the other
# "ej påbörjat" rows. The selection cursor also moved to A16 afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 - this shifts rows 6:19 down to 7:20,
# carrying their content/styles/row-heights along (matches the diff, where
# every row from the old row 6 onward reappears one row lower, unchanged).
$ws.Rows.Item(6).Insert()

# Fill in the newly inserted row 6 with the new task and its status.
$ws.Range("A6").Value = "Distance sampling-analys på ripdata"
$ws.Range("B6").Value = "ej påbörjat"

# Match the red font styling used by the other "ej påbörjat" rows (e.g. the
# row that is now A14/B14, "GIS-data lyornas avstånd till rödrävslyor").
$ws.Range("B6").Font.Color = $ws.Range("B14").Font.Color

# The saved file's selection moved to A16.
$ws.Range("A16").Select()
